# Update imputed KNN values for several cells in the B and D columns.
# (commit message: "Update Name of Algo" -- rerun of the imputation
# algorithm produced slightly different numeric results for these cells)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -7.285000000000001
$ws.Range("D4").Value = -7.595000000000001
$ws.Range("B11").Value = 6.598000000000001
$ws.Range("B12").Value = 5.325
$ws.Range("D14").Value = -7.390000000000001
$ws.Range("B15").Value = 5.538
$ws.Range("D26").Value = -8.119
$ws.Range("B27").Value = 6.531999999999999
$ws.Range("B28").Value = 6.725999999999999
$ws.Range("B31").Value = 6.647
$ws.Range("D31").Value = -8.362
$ws.Range("B32").Value = 7.341999999999999
$ws.Range("D35").Value = -8.013000000000002
$ws.Range("B36").Value = 8.548
$ws.Range("D37").Value = -7.793000000000001
$ws.Range("B38").Value = 6.052000000000001
$ws.Range("D39").Value = -7.397999999999999
$ws.Range("D40").Value = -7.976000000000001
$ws.Range("D45").Value = -7.676
$ws.Range("B46").Value = 6.477000000000001
$ws.Range("D52").Value = -7.961000000000001
$ws.Range("B54").Value = 5.050000000000001
$ws.Range("B55").Value = 4.791
$ws.Range("B56").Value = 4.843999999999999
$ws.Range("D57").Value = -8.100999999999999
$ws.Range("B67").Value = 4.915999999999999
$ws.Range("B69").Value = 5.063000000000001
$ws.Range("B72").Value = 5.684
$ws.Range("B73").Value = 6.601000000000001
$ws.Range("D81").Value = -6.877
$ws.Range("B83").Value = 5.633
$ws.Range("D83").Value = -8.572000000000001
$ws.Range("B86").Value = 5.146000000000001
$ws.Range("B91").Value = 6.109
$ws.Range("B93").Value = 5.33
$ws.Range("B99").Value = 5.217
$ws.Range("D100").Value = -8.105
$ws.Range("D102").Value = -7.894
